# Driver 1 reflection and updated algorithm
$d = $word.ActiveDocument

# Fill in the info table: Driver Full Name, Partner Full Name, Student ID
$table = $d.Tables.Item(1)
$table.Cell(1, 2).Range.InsertAfter("Caitlin Burns")
$table.Cell(2, 2).Range.InsertAfter("Leif Labianca")
$table.Cell(3, 2).Range.InsertAfter("1903398")

# Add the reflection paragraph at the end of the document body (after the
# trailing empty paragraph, before the sectPr), starting with a tab.
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()

$tab = [char]9
$reflectionText = $tab + "The objective of this lab was to get experience using python’s math module while also getting experience using calculations in code. Another objective of this lab was to get familiar with using flowcharts and making sure the flow chart and the algorithm are consistent with each other. The steps we used during this lab were doing the test cases, creating the algorithm, creating the flow chart, and finally, coding. The key concepts explored were the math modules and the calculations. The results did match what we expected to get from this lab, and we used a couple of different scenarios from the test cases. The main challenge I encountered was working on the flow chart, because it had to be consistent with the algorithm. I fixed that by going through both of them together. I do think I learned what I was supposed to after doing this lab because now I am more familiar with math and calculations through python. I liked working with my partner because it is interesting to work with different people instead of the same person all of the time. "

$newRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$newRange.InsertAfter($reflectionText)
